$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Polish" header label in C1 was renamed to "Polski"
$ws.Range("C1").Value = "Polski"

# Reflect the user's final selection on the sheet (was C1 after the edit)
$ws.Range("C1").Select()
